# Applies the coin-price refresh captured in the commit diff.
# (prices/volumes updated; a few coin rows also got re-ordered)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking text (e.g. "1.005", "27.133.23")
# that must stay literal text. Flip it to a Text format first so the
# Range.Value assignments below do not get reinterpreted as numbers (which
# would silently drop things like trailing zeros), then flip the style back
# to Normal afterwards so no stray number-format style sticks to the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.133.23'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '1.832.98'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '312.62'
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').Value = '1.005'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').Value = '0.4652'
$ws.Range('E7').Value = '  -0.80%  '
$ws.Range('D8').Value = '0.3714'
$ws.Range('E8').Value = '  +1.75%  '
$ws.Range('D9').Value = '0.07371'
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').Value = '0.8756'
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('D11').Value = '20.03'
$ws.Range('E11').Value = '  -1.01%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.906.67'
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.07797'
$ws.Range('E13').Value = '  +3.60%  '
$ws.Range('D14').Value = '5.352'
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '6.576'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '92.02'
$ws.Range('E16').Value = '  -0.89%  '
$ws.Range('D17').Value = '1.006'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = '0.000008853'
$ws.Range('E18').Value = '  +1.55%  '
$ws.Range('D19').Value = '1.003'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '14.71'
$ws.Range('E20').Value = '  +0.71%  '
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = '26.807.51'
$ws.Range('E21').Value = '  -2.35%  '
$ws.Range('D22').Value = '5.150'
$ws.Range('E22').Value = '  -1.63%  '
$ws.Range('D23').Value = '10.60'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').Value = '2.139.70'
$ws.Range('E24').Value = '  +3.06%  '
$ws.Range('D25').Value = '152.30'
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('D26').Value = '1.830'
$ws.Range('E26').Value = '  -2.52%  '
$ws.Range('D27').Value = '18.37'
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('D28').Value = '2.098'
$ws.Range('E28').Value = '  -1.70%  '
$ws.Range('D29').Value = '5.092'
$ws.Range('E29').Value = '  -1.43%  '
$ws.Range('D30').Value = '115.89'
$ws.Range('E30').Value = '  -0.40%  '
$ws.Range('D31').Value = '0.08873'
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('D32').Value = '2.959'
$ws.Range('E32').Value = '  +0.59%  '
$ws.Range('D33').Value = '0.7306'
$ws.Range('E33').Value = '  -1.84%  '
$ws.Range('D34').Value = '4.454'
$ws.Range('E34').Value = '  -1.12%  '
$ws.Range('D35').Value = '1.141'
$ws.Range('E35').Value = '  -1.77%  '
$ws.Range('D36').Value = '2.500'
$ws.Range('E36').Value = '  -1.24%  '
$ws.Range('D37').Value = '0.01953'
$ws.Range('E37').Value = '  +1.15%  '
$ws.Range('D38').Value = '1.071'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('D39').Value = '0.05229'
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').Value = '7.134'
$ws.Range('E41').Value = '  -2.45%  '
$ws.Range('D42').Value = '0.5220'
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('D43').Value = '0.8820'
$ws.Range('E43').Value = '  -12.40%  '
$ws.Range('D44').Value = '0.1632'
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('D45').Value = '8.242'
$ws.Range('E45').Value = '  -1.48%  '
$ws.Range('D46').Value = '0.4839'
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('D47').Value = '10.26'
$ws.Range('E47').Value = '  -1.50%  '
$ws.Range('D48').Value = '1.005'
$ws.Range('E48').Value = '  -0.15%  '
$ws.Range('D49').Value = '102.45'
$ws.Range('E49').Value = '  -1.68%  '
$ws.Range('D50').Value = '1.627'
$ws.Range('E50').Value = '  -1.30%  '
$ws.Range('D51').Value = '0.06229'
$ws.Range('E51').Value = '  -0.54%  '

# Restore the default style on column D now that the text is committed.
$priceRange.Style = "Normal"
